# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy the formatting used by the rest of the
# header row (bold, bordered, centered) from A1, then set their labels.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in team record data (Wins=92, Losses=70, Ties=0) for every data row.
$ws.Range("AD2:AD37").Value = 92
$ws.Range("AE2:AE37").Value = 70
$ws.Range("AF2:AF37").Value = 0
